# Update Name of Algo - update column C values for specific rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = -13.021
$ws.Range("C12").Value = -10.616
$ws.Range("C18").Value = -12.187
$ws.Range("C37").Value = -13.075
$ws.Range("C55").Value = -13.837
$ws.Range("C68").Value = -10.809
$ws.Range("C77").Value = -12.954
$ws.Range("C78").Value = -12.798
